$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header row (row 1): field names
$ws.Range("G1").Value = "simple_atlas"
$ws.Range("H1").Value = "simple_tex"

# Row 2: types
$ws.Range("G2").Value = "string"
$ws.Range("H2").Value = "string"

# Row 3: Chinese descriptions
$ws.Range("G3").Value = "缩略图集"
$ws.Range("H3").Value = "缩略图"

# Row 4: data value (only column G gets a value)
$ws.Range("G4").Value = "CardSimple"

# Update the selected cell to match the saved view state
$ws.Range("I4").Select()
